$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86, shifting existing rows 86-110 down to 87-111.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new week's data.
$ws.Cells.Item(86, 1).Value = 8
$ws.Cells.Item(86, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = 44951
$ws.Cells.Item(86, 5).Value = 4
$ws.Cells.Item(86, 6).Value = 100112030
$ws.Cells.Item(86, 7).Value = "Poroto granado"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 400
$ws.Cells.Item(86, 11).Value = 37500
$ws.Cells.Item(86, 12).Value = 38000
$ws.Cells.Item(86, 13).Value = 37750
$ws.Cells.Item(86, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(86, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(86, 16).Value = 1510
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = "Hortaliza"
